# Update 想去人数 (F column) counts across all four sheets.
# This mirrors a re-scrape of the source site where each event's
# "want to go" counter ticked up by a small amount.

$wb = $excel.ActiveWorkbook

# 展览 (sheet 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 352
$ws1.Range("F3").Value  = 3547
$ws1.Range("F4").Value  = 375
$ws1.Range("F7").Value  = 113
$ws1.Range("F8").Value  = 2206
$ws1.Range("F9").Value  = 17
$ws1.Range("F13").Value = 626
$ws1.Range("F20").Value = 57198
$ws1.Range("F30").Value = 4971
$ws1.Range("F35").Value = 1294
$ws1.Range("F36").Value = 1574
$ws1.Range("F42").Value = 36
$ws1.Range("F44").Value = 232

# 演出 (sheet 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 118
$ws2.Range("F44").Value = 24

# 本地生活 (sheet 3)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value  = 1581
$ws3.Range("F9").Value  = 9401
$ws3.Range("F16").Value = 2169
$ws3.Range("F18").Value = 454

# 全部类型 (sheet 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 352
$ws4.Range("F3").Value  = 3547
$ws4.Range("F4").Value  = 375
$ws4.Range("F5").Value  = 1581
$ws4.Range("F11").Value = 113
$ws4.Range("F13").Value = 2169
$ws4.Range("F15").Value = 626
$ws4.Range("F18").Value = 57198
$ws4.Range("F26").Value = 4971
$ws4.Range("F30").Value = 1294
$ws4.Range("F31").Value = 118
$ws4.Range("F32").Value = 454
$ws4.Range("F48").Value = 24
